$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.617.54"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.451.67"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'578.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'145.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("D7").Value = "3.452.82"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").Value = "'7.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").Value = "'0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "4.042.66"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  +7.42%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").Value = "3.454.41"
$ws.Range("D18").Value = "61.735.76"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "'6.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.00%  "
$ws.Range("D20").Value = "'14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.22%  "
$ws.Range("D21").Value = "'9.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").Value = "'403.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("D23").Value = "'0.568"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "'74.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").Value = "'5.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").Value = "3.592.38"
$ws.Range("E29").Value = "  +3.67%  "
$ws.Range("D30").Value = "'7.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'8.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  -9.70%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'23.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").Value = "'7.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "3.478.26"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'5.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "'166.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'0.0792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").Value = "'27.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("D44").Value = "'0.804"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("D45").Value = "'4.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'42.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "2.614.00"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'1.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  +2.64%  "
